$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 35-37 hold the daily scores for 2025-02-12.
# Column A holds dates formatted as plain text (e.g. "2025-02-01") in the
# existing data, so force a text number format before assigning the value
# to stop the engine from auto-converting the string into a date serial,
# then reset the style back to Normal so no stray formatting is left
# behind on the cell.
$ws.Range("A35:A37").NumberFormat = "@"

$ws.Range("A35").Value = "2025-02-12"
$ws.Range("B35").Value = "sleep"
$ws.Range("C35").Value = $true
$ws.Range("D35").Value = $true

$ws.Range("A36").Value = "2025-02-12"
$ws.Range("B36").Value = "activity"
$ws.Range("C36").Value = $true
$ws.Range("D36").Value = $true

$ws.Range("A37").Value = "2025-02-12"
$ws.Range("B37").Value = "weekly_activity"
$ws.Range("C37").Value = $false
$ws.Range("D37").Value = $false

$ws.Range("A35:A37").Style = "Normal"
